# "Fixed error in excel example"
#
# The sheet's siren column (A2:A3) had been imported from the source CSV
# as text (shared string "110 043 015") instead of a real number. This
# fixes it by writing the actual numeric value, which also drops the
# now-unused shared string and the workbook-level defined name that
# pointed at the (external-data-import) range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the workbook-scoped defined name "exemple_valide" (Feuil1!$A$1:$H$3)
if ($wb.Names.Count -gt 0) {
    for ($i = $wb.Names.Count; $i -ge 1; $i--) {
        $wb.Names.Item($i).Delete()
    }
}

# siren was stored as text "110 043 015"; correct it to the numeric value
$ws.Range("A2").Value = 110043015
$ws.Range("A3").Value = 110043015

# Restore the saved selection/cursor position
$ws.Range("B6").Select() | Out-Null
